# Update values in result_data_KNN.xlsx to reflect refreshed KNN imputation
# output ("Update Name of Algo" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = -8.233999999999998
    "D21" = -8.1
    "D23" = -7.874
    "D25" = -8.052000000000001
    "E27" = 16.727
    "E31" = 16.675
    "E39" = 16.604
    "E48" = 17.062
    "E51" = 16.617
    "E52" = 16.543
    "D53" = -7.755
    "E55" = 16.416
    "E56" = 16.276
    "D57" = -7.923999999999999
    "E57" = 16.568
    "D59" = -8.061
    "D69" = -7.597
    "E73" = 16.572
    "D79" = -7.772
    "D83" = -8.237
    "E89" = 17.362
    "E90" = 16.73
    "D93" = -7.516999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
